# Bugfix/Update wrong column names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mistyped "Provider _ID" header (space before ID) -> "Provider__ID"
$ws.Range("T1").Value = "Provider__ID"

# Reset view/selection to the top-left of the sheet
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C11").Select()
